$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '51.791.48'
$ws.Range('E2').Value = '  +1.56%  '
$ws.Range('D3').Value = '2.807.57'
$ws.Range('E3').Value = '  +1.97%  '
$ws.Range('D5').Value = '''353.09'
$ws.Range('E5').Value = '  +0.53%  '
$ws.Range('D6').Value = '''112.30'
$ws.Range('E6').Value = '  +4.77%  '
$ws.Range('D7').Value = '''0.559'
$ws.Range('E7').Value = '  +2.51%  '
$ws.Range('E8').Value = '  +0.03%  '
$ws.Range('D9').Value = '''0.623'
$ws.Range('E9').Value = '  +7.70%  '
$ws.Range('D10').Value = '''40.23'
$ws.Range('E10').Value = '  +2.99%  '
$ws.Range('E11').Value = '  -0.51%  '
$ws.Range('D12').Value = '''0.0839'
$ws.Range('E12').Value = '  +0.92%  '
$ws.Range('D13').Value = '''19.92'
$ws.Range('E13').Value = '  +1.80%  '
$ws.Range('D14').Value = '''7.78'
$ws.Range('E14').Value = '  +4.46%  '
$ws.Range('D15').Value = '3.243.46'
$ws.Range('E15').Value = '  +1.90%  '
$ws.Range('D16').Value = '2.805.22'
$ws.Range('E16').Value = '  +1.93%  '
$ws.Range('D17').Value = '''0.954'
$ws.Range('E17').Value = '  +3.34%  '
$ws.Range('D18').Value = '51.806.86'
$ws.Range('E18').Value = '  +1.81%  '
$ws.Range('D19').Value = '''7.63'
$ws.Range('E19').Value = '  +0.72%  '
$ws.Range('D20').Value = '''3.31'
$ws.Range('E20').Value = '  +9.28%  '
$ws.Range('E21').Value = '  +4.53%  '
$ws.Range('D22').Value = '0.0₃0974'
$ws.Range('E22').Value = '  +2.21%  '
$ws.Range('D23').Value = '''70.30'
$ws.Range('E23').Value = '  +1.71%  '
$ws.Range('D24').Value = '''267.80'
$ws.Range('E24').Value = '  +1.90%  '
$ws.Range('E25').Value = '  +2.41%  '
$ws.Range('E26').Value = '  +0.05%  '
$ws.Range('D27').Value = '''26.15'
$ws.Range('E27').Value = '  +1.40%  '
$ws.Range('D28').Value = '''0.162'
$ws.Range('E28').Value = '  +0.71%  '
$ws.Range('D29').Value = '''39.32'
$ws.Range('E29').Value = '  +14.80%  '
$ws.Range('D30').Value = '''10.42'
$ws.Range('E30').Value = '  +4.11%  '
$ws.Range('E31').Value = '  +2.13%  '
$ws.Range('D32').Value = '''52.49'
$ws.Range('E32').Value = '  +2.10%  '
$ws.Range('D33').Value = '''6.12'
$ws.Range('E33').Value = '  +1.98%  '
$ws.Range('D34').Value = '''0.0906'
$ws.Range('E34').Value = '  +9.57%  '
$ws.Range('D35').Value = '''0.0450'
$ws.Range('E35').Value = '  +2.58%  '
$ws.Range('E36').Value = '  +4.59%  '
$ws.Range('D37').Value = '''0.999'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').Value = '''19.07'
$ws.Range('E38').Value = '  +4.73%  '
$ws.Range('B39').Value = 'ARBITRUM'
$ws.Range('C39').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D39').Value = '''2.02'
$ws.Range('E39').Value = '  +3.82%  '
$ws.Range('B40').Value = 'LidoDAOToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D40').Value = '''3.17'
$ws.Range('E40').Value = '  +1.58%  '
$ws.Range('E41').Value = '  +2.45%  '
$ws.Range('E42').Value = '  +1.78%  '
$ws.Range('E43').Value = '  +2.04%  '
$ws.Range('D44').Value = '''120.21'
$ws.Range('E44').Value = '  -0.41%  '
$ws.Range('D45').Value = '''21.82'
$ws.Range('E45').Value = '  -0.70%  '
$ws.Range('D46').Value = '''3.55'
$ws.Range('E46').Value = '  +10.80%  '
$ws.Range('E47').Value = '  +9.67%  '
$ws.Range('D48').Value = '2.120.34'
$ws.Range('E48').Value = '  +1.68%  '
$ws.Range('E49').Value = '  +7.65%  '
$ws.Range('B50').Value = 'TrustWalletToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D50').Value = '''1.38'
$ws.Range('E50').Value = '  +8.41%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '''5.49'
$ws.Range('E51').Value = '  +1.43%  '
